{"js": "const NEW_TEXT_XML = \"    &lt;---asTable(java.lang.String,java.lang.String,java.lang.String,java.lang.String,java.lang.String) with arguments [excel.xlsx, NotExistingSheet, C3, F7, fr-FR] failed:\\n\\tjava.lang.IllegalArgumentException: The sheet NotExistingSheet doesn't exist in file:/home/development/git/M2Doc/tests/org.obeonetwork.m2doc.tests/resources/excelServices/asTableNotExistingSheet/excel.xlsx.\\n\\t\\tat org.obeonetwork.m2doc.services.ExcelServices.asTable(ExcelServices.java:129)\\n\\t\\tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)\\n\\t\\tat java.base/java.lang.reflect.Method.invoke(Method.java:580)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:170)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:231)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallService(EvaluationServices.java:122)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCall(EvaluationServices.java:237)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallOrApply(EvaluationServices.java:273)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:173)\\n\\t\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:244)\\n\\t\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:146)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:53)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:674)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)\\n\\t\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\\n\\t\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:350)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\\n\\t\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:334)\\n\\t\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:878)\\n\\t\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:499)\\n\\t\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)\\n\\t\\tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)\\n\\t\\tat java.base/java.lang.reflect.Method.invoke(Method.java:580)\\n\\t\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:59)\\n\\t\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\\n\\t\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)\\n\\t\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\\n\\t\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\\n\\t\\tat org.junit.runners.BlockJUnit4ClassRunner$1.evaluate(BlockJUnit4ClassRunner.java:100)\\n\\t\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)\\n\\t\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)\\n\\t\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)\\n\\t\\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\\n\\t\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\\n\\t\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\\n\\t\\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\\n\\t\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\\n\\t\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\\n\\t\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\\n\\t\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\\n\\t\\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\\n\\t\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\\n\\t\\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\\n\\t\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\\n\\t\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\\n\\t\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\\n\\t\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\\n\\t\\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\\n\\t\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\\n\\t\\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\\n\\t\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\\n\\t\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)\\n\\t\";\n\nconst results = context.document.body.search(\"    <---asTable\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor text not found\");\n}\n\nconst anchor = results.items[0];\nconst para = anchor.paragraphs.getFirst();\nconst paraEnd = para.getRange(\"End\");\nconst target = anchor.expandTo(paraEnd);\n\nconst runXml =\n  '<w:r><w:rPr><w:b w:val=\"on\"/><w:color w:val=\"FF0000\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">' + NEW_TEXT_XML + '</w:t><w:br/></w:r>';\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' + runXml + '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"    <---asTable\")\nif (-not $found) {\n    throw \"Anchor text not found\"\n}\n$matchStart = $rng.Start\n\n$para = $d.Paragraphs.Item(2)\n$targetEnd = $para.Range.End - 2\n\n$newRng = $d.Range($matchStart, $targetEnd)\n\n$newText = @'\n    <---asTable(java.lang.String,java.lang.String,java.lang.String,java.lang.String,java.lang.String) with arguments [excel.xlsx, NotExistingSheet, C3, F7, fr-FR] failed:\n\tjava.lang.IllegalArgumentException: The sheet NotExistingSheet doesn't exist in file:/home/development/git/M2Doc/tests/org.obeonetwork.m2doc.tests/resources/excelServices/asTableNotExistingSheet/excel.xlsx.\n\t\tat org.obeonetwork.m2doc.services.ExcelServices.asTable(ExcelServices.java:129)\n\t\tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)\n\t\tat java.base/java.lang.reflect.Method.invoke(Method.java:580)\n\t\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:170)\n\t\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:231)\n\t\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallService(EvaluationServices.java:122)\n\t\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCall(EvaluationServices.java:237)\n\t\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallOrApply(EvaluationServices.java:273)\n\t\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:173)\n\t\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:244)\n\t\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\t\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:146)\n\t\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:53)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:674)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)\n\t\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\n\t\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:350)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\n\t\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\t\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\n\t\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:334)\n\t\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:878)\n\t\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:499)\n\t\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)\n\t\tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)\n\t\tat java.base/java.lang.reflect.Method.invoke(Method.java:580)\n\t\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:59)\n\t\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\n\t\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)\n\t\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\n\t\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\t\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\n\t\tat org.junit.runners.BlockJUnit4ClassRunner$1.evaluate(BlockJUnit4ClassRunner.java:100)\n\t\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)\n\t\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)\n\t\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)\n\t\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\n\t\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\n\t\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\n\t\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\n\t\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\n\t\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\n\t\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\t\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\t\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\n\t\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\n\t\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\n\t\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\n\t\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\n\t\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\n\t\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\t\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\n\t\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\n\t\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\t\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\t\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\n\t\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\n\t\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\n\t\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\n\t\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\n\t\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\n\t\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\n\t\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)\n\t\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)\n\t\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)\n\t\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)\n\t\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)\n\t\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)\n\t\n'@\n\n$newRng.Text = $newText\n"}
